# bore_data.xlsx - "updates to try get working structural model"
#
# - Row 16 on the "geo" sheet gets reclassified the same way several other
#   bores already were: it becomes a "Control" / "Top unit(s) pinched out"
#   row (same style + text pattern as e.g. row 19), its Kp thickness (G) is
#   recorded as 0, and the stale Jy thickness (M) is cleared.
# - A handful of other already-"Control" rows get their Kp thickness (G)
#   filled in with 0 instead of being left blank.
# - Two more rows (M26, M53) have their stale Jy thickness cleared, matching
#   the pattern used elsewhere once a bore becomes a Control point.
# - The "geo" sheet becomes the active/selected sheet (whole used columns
#   A:S selected), and "strat" loses the tab-selected flag it used to have.

$wb = $excel.ActiveWorkbook
$geo = $wb.Worksheets.Item("geo")
$strat = $wb.Worksheets.Item("strat")

# --- Row 16: turn into a "Control" / pinched-out row, matching the format
#     already used by rows like 7, 19, 26, ... (copy formatting from row 19,
#     a row with an identical starting layout, then set the new values).
$geo.Range("A19:G19").Copy()
$geo.Range("A16:G16").PasteSpecial(-4122)  # xlPasteFormats
$geo.Range("D16").Value = "Control"
$geo.Range("E16").Value = "Top unit(s) pinched out"
$geo.Range("G16").Value = 0
$geo.Range("M16").ClearContents()

# --- Kp thickness (G) recorded as 0 for these already-"Control" rows.
$gZeroRows = 7,19,26,29,31,35,39,45,47,49,51,53,56
foreach ($r in $gZeroRows) {
    $geo.Cells.Item($r, 7).Value = 0
}

# --- Stale Jy thickness (M) cleared for these rows.
$mClearRows = 26,53
foreach ($r in $mClearRows) {
    $geo.Cells.Item($r, 13).ClearContents()
}

# --- View state: make "geo" the active sheet/tab and select the full used
#     range, which also drops "strat"'s tabSelected flag.
$geo.Activate()
$geo.Range("A1:S1048576").Select()
